$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "-"
$ws.Range("C3").Value = "-"
$ws.Range("E3").Value = "Andre Lucca-Circuitos Elétricos"
$ws.Range("C6").Value = "João Rodrigues-Desenho Técnico"
$ws.Range("E6").Value = "-"
$ws.Range("C7").Value = "João Rodrigues-Desenho Técnico"
